$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I1 ("I0") and J1 ("IF"), styled like the other headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats - reuse the existing header style

# Data rows 2-9: I column is always 1, J column mirrors the H column value
for ($r = 2; $r -le 9; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
